$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "What was the production of oil in the state of Rio de Janeiro?"
$ws.Range("B2").Value = "What was the production of oil in the state of Rio de Janeiro?"

# Row 3
$ws.Range("A3").Value = "What was the average monthly production of oil in the state of Rio de Janeiro?"
$ws.Range("B3").Value = "What was the month production of oil in the state of Rio de Janeiro?"

# Row 4
$ws.Range("A4").Value = "What was the average yearly production of oil in the state of Alagoas?"
$ws.Range("B4").Value = "What was the year production of oil in the state of Alagoas?"

# Row 6
$ws.Range("A6").Value = "What was the maximum production of oil in the state of Ceará per field?"
$ws.Range("B6").Value = "What was the production of oil in the state of Ceará field?"

# Row 7
$ws.Range("A7").Value = "What was the minimum gas production in the state of São Paulo per basin?"
$ws.Range("B7").Value = "What was the gas production in the state of São Paulo basin?"

# Row 8
$ws.Range("A8").Value = "What was the average monthly oil production by the operator Petrobrás?"
$ws.Range("B8").Value = "What was the month oil production the operator Petrobrás?"
$ws.Range("C8").Value = "SELECT year, month, oil_production, operator FROM ANP WHERE lower(operator) = 'petrobras'"
$ws.Range("D8").Value = "select operator, avg(oil_production) as avg_oil_production from (select operator, sum(oil_production) as oil_production from nlidb_result_set group by operator, year, month) group by operator order by operator"

# Row 13
$ws.Range("A13").Value = "What was the mean monthly petroleum production by field in the state of Rio de Janeiro?"
$ws.Range("B13").Value = "What was the month petroleum production field in the state of Rio de Janeiro?"

# Row 15
$ws.Range("A15").Value = "What was the average monthly production of oil per field in the state of Rio de Janeiro and year 2015?"
$ws.Range("B15").Value = "What was the month production of oil field in the state of Rio de Janeiro year 2015?"
$ws.Range("D15").Value = "select field, year, avg(oil_production) as avg_oil_production from (select field, year, sum(oil_production) as oil_production from nlidb_result_set group by field, year, year, month) group by field, year order by field, year"

# Row 16
$ws.Range("A16").Value = "What was the average yearly production of oil per field and state in the year in 2015?"
$ws.Range("B16").Value = "What was the year production of oil field state in the year in 2015?"

$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.Range("D16").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
